$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(57, 8).Value = 10000
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 10000
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 30000
$ws.Cells.Item(57, 14).Value = -30998

$ws.Cells.Item(113, 8).Value = 17514.445
$ws.Cells.Item(113, 9).Value = 7411.3335
$ws.Cells.Item(113, 10).Value = 22566
$ws.Cells.Item(113, 11).Value = 7411.3335
$ws.Cells.Item(113, 12).Value = 22566
$ws.Cells.Item(113, 13).Value = -4157.3335
$ws.Cells.Item(113, 14).Value = -29074

$ws.Cells.Item(132, 8).Value = 4168.836
$ws.Cells.Item(132, 9).Value = 2628.585
$ws.Cells.Item(132, 10).Value = 14373
$ws.Cells.Item(132, 11).Value = 7885.755
$ws.Cells.Item(132, 12).Value = 43119
$ws.Cells.Item(132, 13).Value = -5355.755
$ws.Cells.Item(132, 14).Value = -48179

$ws.Cells.Item(136, 8).Value = 86744.164
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 86744.164
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 86744.164
$ws.Cells.Item(136, 14).Value = -96944.164

$ws.Cells.Item(137, 8).Value = 2281.2334
$ws.Cells.Item(137, 9).Value = 2456.56
$ws.Cells.Item(137, 10).Value = 2156
$ws.Cells.Item(137, 11).Value = 7369.68
$ws.Cells.Item(137, 12).Value = 6468
$ws.Cells.Item(137, 13).Value = -4819.68
$ws.Cells.Item(137, 14).Value = -11568

$ws.Cells.Item(138, 8).Value = 3273.4927
$ws.Cells.Item(138, 9).Value = 2805.8125
$ws.Cells.Item(138, 10).Value = 3414.6792
$ws.Cells.Item(138, 11).Value = 8417.4375
$ws.Cells.Item(138, 12).Value = 10244.0376
$ws.Cells.Item(138, 13).Value = -3277.4375
$ws.Cells.Item(138, 14).Value = -20524.0376


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6915.531
$ws.Cells.Item(32, 9).Value = 4105.986
$ws.Cells.Item(32, 10).Value = 24794.455
$ws.Cells.Item(32, 11).Value = 4105.986
$ws.Cells.Item(32, 12).Value = 24794.455
$ws.Cells.Item(32, 13).Value = -3818.986
$ws.Cells.Item(32, 14).Value = -25368.455

$ws.Cells.Item(45, 8).Value = 7245.4067
$ws.Cells.Item(45, 9).Value = 14025.263
$ws.Cells.Item(45, 10).Value = 4024.975
$ws.Cells.Item(45, 11).Value = 14025.263
$ws.Cells.Item(45, 12).Value = 4024.975
$ws.Cells.Item(45, 13).Value = -13648.263
$ws.Cells.Item(45, 14).Value = -4778.975

$ws.Cells.Item(92, 8).Value = 35366.332
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 35366.332
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 35366.332
$ws.Cells.Item(92, 14).Value = -40358.332

$ws.Cells.Item(102, 8).Value = 9668.895
$ws.Cells.Item(102, 9).Value = 9928.277
$ws.Cells.Item(102, 10).Value = 5000
$ws.Cells.Item(102, 11).Value = 9928.277
$ws.Cells.Item(102, 12).Value = 5000
$ws.Cells.Item(102, 13).Value = -8306.277
$ws.Cells.Item(102, 14).Value = -8244

$ws.Cells.Item(110, 8).Value = 3004
$ws.Cells.Item(110, 9).Value = 2645.647
$ws.Cells.Item(110, 10).Value = 6050
$ws.Cells.Item(110, 11).Value = 2645.647
$ws.Cells.Item(110, 12).Value = 6050
$ws.Cells.Item(110, 13).Value = -600.6469999999999
$ws.Cells.Item(110, 14).Value = -10140

$ws.Cells.Item(129, 8).Value = 65500
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 65500
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 65500
$ws.Cells.Item(129, 14).Value = -75500


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(57, 8).Value = 88121.86
$ws.Cells.Item(57, 9).Value = 77709
$ws.Cells.Item(57, 10).Value = 88922.84
$ws.Cells.Item(57, 11).Value = 77709
$ws.Cells.Item(57, 12).Value = 88922.84
$ws.Cells.Item(57, 13).Value = -76989
$ws.Cells.Item(57, 14).Value = -90362.84

$ws.Cells.Item(60, 8).Value = 28781.334
$ws.Cells.Item(60, 9).Value = 22709
$ws.Cells.Item(60, 10).Value = 29995.8
$ws.Cells.Item(60, 11).Value = 22709
$ws.Cells.Item(60, 12).Value = 29995.8
$ws.Cells.Item(60, 13).Value = -22110
$ws.Cells.Item(60, 14).Value = -31193.8

$ws.Cells.Item(100, 8).Value = 36536
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 36536
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 36536
$ws.Cells.Item(100, 14).Value = -38700

$ws.Cells.Item(136, 8).Value = 88121.86
$ws.Cells.Item(136, 9).Value = 77709
$ws.Cells.Item(136, 10).Value = 88922.84
$ws.Cells.Item(136, 11).Value = 77709
$ws.Cells.Item(136, 12).Value = 88922.84
$ws.Cells.Item(136, 13).Value = -72609
$ws.Cells.Item(136, 14).Value = -99122.84

$ws.Cells.Item(139, 8).Value = 89155.37
$ws.Cells.Item(139, 9).Value = 80709
$ws.Cells.Item(139, 10).Value = 90000
$ws.Cells.Item(139, 11).Value = 80709
$ws.Cells.Item(139, 12).Value = 90000
$ws.Cells.Item(139, 13).Value = -75569
$ws.Cells.Item(139, 14).Value = -100280

$ws.Cells.Item(141, 8).Value = 89890
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 89890
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 89890
$ws.Cells.Item(141, 14).Value = -100250


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7834.85
$ws.Cells.Item(31, 9).Value = 4877.1924
$ws.Cells.Item(31, 10).Value = 13327.643
$ws.Cells.Item(31, 11).Value = 4877.1924
$ws.Cells.Item(31, 12).Value = 13327.643
$ws.Cells.Item(31, 13).Value = -4582.1924
$ws.Cells.Item(31, 14).Value = -13917.643

$ws.Cells.Item(34, 8).Value = 7834.85
$ws.Cells.Item(34, 9).Value = 4877.1924
$ws.Cells.Item(34, 10).Value = 13327.643
$ws.Cells.Item(34, 11).Value = 4877.1924
$ws.Cells.Item(34, 12).Value = 13327.643
$ws.Cells.Item(34, 13).Value = -4675.1924
$ws.Cells.Item(34, 14).Value = -13731.643

$ws.Cells.Item(43, 8).Value = 30000
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 30000
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 30000
$ws.Cells.Item(43, 14).Value = -30368

$ws.Cells.Item(101, 8).Value = 30000
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 30000
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 30000
$ws.Cells.Item(101, 14).Value = -36490

$ws.Cells.Item(138, 8).Value = 88655.57000000001
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 88655.57000000001
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 88655.57000000001
$ws.Cells.Item(138, 14).Value = -98935.57000000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 5300
$ws.Cells.Item(3, 9).Value = 3000
$ws.Cells.Item(3, 10).Value = 9900
$ws.Cells.Item(3, 11).Value = 9000
$ws.Cells.Item(3, 12).Value = 29700
$ws.Cells.Item(3, 13).Value = -8888
$ws.Cells.Item(3, 14).Value = -29924

$ws.Cells.Item(8, 8).Value = 127533.57
$ws.Cells.Item(8, 9).Value = 127533.57
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 382600.71
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = -382461.71

$ws.Cells.Item(46, 8).Value = 22987.617
$ws.Cells.Item(46, 9).Value = 1268.8718
$ws.Cells.Item(46, 10).Value = 128866.5
$ws.Cells.Item(46, 11).Value = 3806.6154
$ws.Cells.Item(46, 12).Value = 386599.5
$ws.Cells.Item(46, 13).Value = -3715.6154
$ws.Cells.Item(46, 14).Value = -386781.5

$ws.Cells.Item(80, 8).Value = 3047.5
$ws.Cells.Item(80, 9).Value = 540
$ws.Cells.Item(80, 10).Value = 5555
$ws.Cells.Item(80, 11).Value = 1620
$ws.Cells.Item(80, 12).Value = 16665
$ws.Cells.Item(80, 13).Value = -684
$ws.Cells.Item(80, 14).Value = -18537

$ws.Cells.Item(83, 8).Value = 3047.5
$ws.Cells.Item(83, 9).Value = 540
$ws.Cells.Item(83, 10).Value = 5555
$ws.Cells.Item(83, 11).Value = 4860
$ws.Cells.Item(83, 12).Value = 49995
$ws.Cells.Item(83, 13).Value = -180
$ws.Cells.Item(83, 14).Value = -59355

$ws.Cells.Item(129, 8).Value = 1582.2307
$ws.Cells.Item(129, 9).Value = 1217.4445
$ws.Cells.Item(129, 10).Value = 2403
$ws.Cells.Item(129, 11).Value = 3652.3335
$ws.Cells.Item(129, 12).Value = 7209
$ws.Cells.Item(129, 13).Value = 1347.6665
$ws.Cells.Item(129, 14).Value = -17209

$ws.Cells.Item(136, 8).Value = 10578.083
$ws.Cells.Item(136, 9).Value = 11159.667
$ws.Cells.Item(136, 10).Value = 8833.333000000001
$ws.Cells.Item(136, 11).Value = 33479.001
$ws.Cells.Item(136, 12).Value = 26499.999
$ws.Cells.Item(136, 13).Value = -28379.001
$ws.Cells.Item(136, 14).Value = -36699.999


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(101, 8).Value = 80833.734
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 80833.734
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 80833.734
$ws.Cells.Item(101, 14).Value = -87323.734

$ws.Cells.Item(122, 8).Value = 3447.182
$ws.Cells.Item(122, 9).Value = 2059.8572
$ws.Cells.Item(122, 10).Value = 5875
$ws.Cells.Item(122, 11).Value = 6179.571599999999
$ws.Cells.Item(122, 12).Value = 17625
$ws.Cells.Item(122, 13).Value = -3729.571599999999
$ws.Cells.Item(122, 14).Value = -22525

$ws.Cells.Item(135, 8).Value = 86510.086
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 86510.086
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 86510.086
$ws.Cells.Item(135, 14).Value = -96650.086

$ws.Cells.Item(139, 8).Value = 89165.57000000001
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 89165.57000000001
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 89165.57000000001
$ws.Cells.Item(139, 14).Value = -99445.57000000001


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1249.7858
$ws.Cells.Item(46, 9).Value = 1419.8
$ws.Cells.Item(46, 10).Value = 1155.3334
$ws.Cells.Item(46, 11).Value = 1419.8
$ws.Cells.Item(46, 12).Value = 1155.3334
$ws.Cells.Item(46, 13).Value = -1231.8
$ws.Cells.Item(46, 14).Value = -1531.3334

$ws.Cells.Item(69, 8).Value = 90000
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 90000
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 90000
$ws.Cells.Item(69, 14).Value = -91622

$ws.Cells.Item(72, 8).Value = 90000
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 90000
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 270000
$ws.Cells.Item(72, 14).Value = -278112

$ws.Cells.Item(101, 8).Value = 48890
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 48890
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 48890
$ws.Cells.Item(101, 14).Value = -55380


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3276.087
$ws.Cells.Item(126, 9).Value = 2118.2
$ws.Cells.Item(126, 10).Value = 10995.333
$ws.Cells.Item(126, 11).Value = 6354.599999999999
$ws.Cells.Item(126, 12).Value = 32985.999
$ws.Cells.Item(126, 13).Value = -3884.599999999999
$ws.Cells.Item(126, 14).Value = -37925.999

$ws.Cells.Item(133, 8).Value = 109750
$ws.Cells.Item(133, 9).Value = 90000
$ws.Cells.Item(133, 10).Value = 116333.336
$ws.Cells.Item(133, 11).Value = 90000
$ws.Cells.Item(133, 12).Value = 116333.336
$ws.Cells.Item(133, 13).Value = -84940
$ws.Cells.Item(133, 14).Value = -126453.336

$ws.Cells.Item(135, 8).Value = 88713.97
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 88713.97
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 88713.97
$ws.Cells.Item(135, 14).Value = -98853.97

